$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 418; existing rows 418:500 shift down to 419:501
$ws.Rows(418).Insert()

# Populate the newly inserted row 418 with its data
$ws.Cells.Item(418, 1).Value = 5
$ws.Cells.Item(418, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(418, 3).Value = "Maule"
$ws.Cells.Item(418, 4).Value = 44995
$ws.Cells.Item(418, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(418, 5).Value = 7
$ws.Cells.Item(418, 6).Value = 100112023
$ws.Cells.Item(418, 7).Value = "Brócoli"
$ws.Cells.Item(418, 8).Value = "Sin especificar"
$ws.Cells.Item(418, 9).Value = "Primera"
$ws.Cells.Item(418, 10).Value = 5000
$ws.Cells.Item(418, 11).Value = 700
$ws.Cells.Item(418, 12).Value = 700
$ws.Cells.Item(418, 13).Value = 700
$ws.Cells.Item(418, 14).Value = "`$/unidad"
$ws.Cells.Item(418, 15).Value = "Región del Maule"
$ws.Cells.Item(418, 16).Value = 700
$ws.Cells.Item(418, 17).Value = 1
$ws.Cells.Item(418, 18).Value = "Hortaliza"
